# Refresh the "cryptos" price list (GitHub Actions scheduled update).
#
# Column D ("Price") holds plain-text figures that often look numeric
# (e.g. "244.80", "44.211.24" — the latter uses '.' as a thousands
# separator, not a decimal point). A bare Range.Value assignment would
# let Excel auto-detect such strings as numbers, silently dropping
# trailing zeros / mangling multi-dot values. Prefixing the literal with
# a leading apostrophe is the standard Excel "force text" entry mode —
# Excel strips the apostrophe and keeps the remainder as literal text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.211.24"
$ws.Range("E2").Value = "  +2.08%  "
$ws.Range("D3").Value = "'2.400.37"
$ws.Range("E3").Value = "  +1.87%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  +8.12%  "
$ws.Range("D6").Value = "'244.80"
$ws.Range("E6").Value = "  +5.06%  "
$ws.Range("D7").Value = "'77.31"
$ws.Range("E7").Value = "  +7.62%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.645"
$ws.Range("E9").Value = "  +34.98%  "
$ws.Range("E10").Value = "  +6.39%  "
$ws.Range("D11").Value = "'57.99"
$ws.Range("E11").Value = "  +2.12%  "
$ws.Range("D12").Value = "'33.90"
$ws.Range("E12").Value = "  +25.54%  "
$ws.Range("D13").Value = "'7.67"
$ws.Range("E13").Value = "  +22.11%  "
$ws.Range("E14").Value = "  +2.41%  "
$ws.Range("D15").Value = "'2.752.53"
$ws.Range("E15").Value = "  +1.84%  "
$ws.Range("D16").Value = "'17.27"
$ws.Range("E16").Value = "  +7.52%  "
$ws.Range("D17").Value = "'0.947"
$ws.Range("E17").Value = "  +9.02%  "
$ws.Range("D18").Value = "'2.399.42"
$ws.Range("E18").Value = "  +1.64%  "
$ws.Range("D19").Value = "'44.220.23"
$ws.Range("E19").Value = "  +2.18%  "
$ws.Range("E20").Value = "  +2.38%  "
$ws.Range("D21").Value = "'6.77"
$ws.Range("E21").Value = "  +6.96%  "
$ws.Range("D22").Value = "'79.40"
$ws.Range("E22").Value = "  +6.92%  "
$ws.Range("D23").Value = "'260.65"
$ws.Range("E23").Value = "  +4.28%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").Value = "'2.57"
$ws.Range("E25").Value = "  +5.27%  "
$ws.Range("D26").Value = "'3.69"
$ws.Range("E26").Value = "  -2.07%  "
$ws.Range("D27").Value = "'11.18"
$ws.Range("E27").Value = "  +11.97%  "
$ws.Range("E28").Value = "  +19.51%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'23.65"
$ws.Range("E29").Value = "  +6.05%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.32"
$ws.Range("E30").Value = "  +2.66%  "
$ws.Range("D31").Value = "'176.05"
$ws.Range("E31").Value = "  +1.45%  "
$ws.Range("E32").Value = "  +1.88%  "
$ws.Range("E33").Value = "  +8.15%  "
$ws.Range("E34").Value = "  +9.36%  "
$ws.Range("D35").Value = "'0.0766"
$ws.Range("E35").Value = "  +11.14%  "
$ws.Range("E36").Value = "  +8.33%  "
$ws.Range("D37").Value = "'3.92"
$ws.Range("E37").Value = "  +6.76%  "
$ws.Range("D38").Value = "'2.53"
$ws.Range("E38").Value = "  +3.81%  "
$ws.Range("E39").Value = "  +1.20%  "
$ws.Range("E40").Value = "  +9.98%  "
$ws.Range("D41").Value = "'19.25"
$ws.Range("E41").Value = "  +3.88%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.204"
$ws.Range("E42").Value = "  +21.40%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'9.12"
$ws.Range("E43").Value = "  +2.33%  "
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("E45").Value = "  +6.55%  "
$ws.Range("E46").Value = "  +16.68%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'105.13"
$ws.Range("E47").Value = "  +5.89%  "
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").Value = "'1.29"
$ws.Range("E48").Value = "  +6.30%  "
$ws.Range("E49").Value = "  +6.88%  "
$ws.Range("D50").Value = "'4.59"
$ws.Range("E50").Value = "  +2.82%  "
$ws.Range("D51").Value = "'55.88"
$ws.Range("E51").Value = "  +11.23%  "
